$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: clear the numeric value -> becomes an empty cell. Nudging the style keeps
# the (now blank) cell present in the sheet instead of being dropped outright.
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = "Normal"

# C4: 16245.68469872195 -> 0
$ws.Range("C4").Value = 0

# C5: 66313.42487543575 -> 0
$ws.Range("C5").Value = 0

# Row 7's label changes from "Other" to "Biogas" and its value is corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 144.8125076895869

# A new row 8 is added, taking over the "Other" label with a new corrected value
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Other"

$ws.Range("B8").Value = ""
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = ""
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = 139.9419385155178
